# "add  your picture and  name to the team section  ." task (row 7) is now
# complete, so mark it "Done" in the "End the task" column, same as the
# other completed rows (C3, C14, C16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C7").Value = "Done"

# Move the active selection (cosmetic, matches the author's final cursor
# position after making the edit).
$ws.Range("E5").Select()
